{"js": "// Update the date line and the 24 changed answer cells in the\n// \"two-digit divided by one-digit\" practice table, in place, while\n// preserving each run's existing character/paragraph formatting.\n\n// 1) Update the date heading paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateRange = paragraphs.items[0].getRange();\ndateRange.insertText(\"2024-09-12 Thursday\", Word.InsertLocation.replace);\n\n// 2) Update the answer cells inside the single practice table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Only the cells that actually change text are listed (row/col are\n// 0-based table coordinates); everything else (fonts, size, alignment,\n// blank spacer rows) is left untouched.\nconst cellUpdates = [\n  { row: 0, col: 0, newText: \"23\u00f78=2, 7\" },\n  { row: 0, col: 1, newText: \"58\u00f75=11, 3\" },\n  { row: 0, col: 2, newText: \"43\u00f78=5, 3\" },\n  { row: 0, col: 3, newText: \"58\u00f72=29, 0\" },\n  { row: 0, col: 4, newText: \"70\u00f75=14, 0\" },\n  { row: 4, col: 0, newText: \"19\u00f78=2, 3\" },\n  { row: 4, col: 1, newText: \"36\u00f79=4, 0\" },\n  { row: 4, col: 2, newText: \"73\u00f76=12, 1\" },\n  { row: 4, col: 3, newText: \"75\u00f75=15, 0\" },\n  { row: 4, col: 4, newText: \"76\u00f73=25, 1\" },\n  { row: 8, col: 0, newText: \"88\u00f77=12, 4\" },\n  { row: 8, col: 1, newText: \"66\u00f78=8, 2\" },\n  { row: 8, col: 2, newText: \"82\u00f76=13, 4\" },\n  { row: 8, col: 3, newText: \"58\u00f74=14, 2\" },\n  { row: 8, col: 4, newText: \"83\u00f73=27, 2\" },\n  { row: 12, col: 0, newText: \"24\u00f76=4, 0\" },\n  { row: 12, col: 1, newText: \"72\u00f79=8, 0\" },\n  { row: 12, col: 2, newText: \"37\u00f74=9, 1\" },\n  { row: 12, col: 3, newText: \"81\u00f79=9, 0\" },\n  { row: 12, col: 4, newText: \"36\u00f78=4, 4\" },\n  { row: 16, col: 1, newText: \"58\u00f74=14, 2\" },\n  { row: 16, col: 2, newText: \"60\u00f75=12, 0\" },\n  { row: 16, col: 3, newText: \"73\u00f79=8, 1\" },\n  { row: 16, col: 4, newText: \"38\u00f78=4, 6\" },\n];\n\nconst cellRanges = cellUpdates.map((u) => {\n  const cell = table.getCell(u.row, u.col);\n  const cellParagraphs = cell.body.paragraphs;\n  cellParagraphs.load(\"items\");\n  return { update: u, cellParagraphs };\n});\n\nawait context.sync();\n\nfor (const { update, cellParagraphs } of cellRanges) {\n  const range = cellParagraphs.items[0].getRange();\n  range.insertText(update.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 24 changed answer cells in the\n# \"two-digit divided by one-digit\" practice table, in place, while\n# preserving each run's existing character/paragraph formatting\n# (setting Range.Text keeps the existing font/size/alignment on the\n# run instead of inserting a brand new default-formatted run).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading paragraph (first paragraph of the body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-09-12 Thursday\"\n\n# 2) Update the answer cells inside the single practice table.\n$t = $d.Tables.Item(1)\n\n# Only the cells that actually change text are listed (Row/Col are\n# 1-based table coordinates, matching the Word object model);\n# everything else (fonts, size, alignment, blank spacer rows) is left\n# untouched.\n$cellUpdates = @(\n  @{ Row = 1; Col = 1; Text = \"23\u00f78=2, 7\" },\n  @{ Row = 1; Col = 2; Text = \"58\u00f75=11, 3\" },\n  @{ Row = 1; Col = 3; Text = \"43\u00f78=5, 3\" },\n  @{ Row = 1; Col = 4; Text = \"58\u00f72=29, 0\" },\n  @{ Row = 1; Col = 5; Text = \"70\u00f75=14, 0\" },\n  @{ Row = 5; Col = 1; Text = \"19\u00f78=2, 3\" },\n  @{ Row = 5; Col = 2; Text = \"36\u00f79=4, 0\" },\n  @{ Row = 5; Col = 3; Text = \"73\u00f76=12, 1\" },\n  @{ Row = 5; Col = 4; Text = \"75\u00f75=15, 0\" },\n  @{ Row = 5; Col = 5; Text = \"76\u00f73=25, 1\" },\n  @{ Row = 9; Col = 1; Text = \"88\u00f77=12, 4\" },\n  @{ Row = 9; Col = 2; Text = \"66\u00f78=8, 2\" },\n  @{ Row = 9; Col = 3; Text = \"82\u00f76=13, 4\" },\n  @{ Row = 9; Col = 4; Text = \"58\u00f74=14, 2\" },\n  @{ Row = 9; Col = 5; Text = \"83\u00f73=27, 2\" },\n  @{ Row = 13; Col = 1; Text = \"24\u00f76=4, 0\" },\n  @{ Row = 13; Col = 2; Text = \"72\u00f79=8, 0\" },\n  @{ Row = 13; Col = 3; Text = \"37\u00f74=9, 1\" },\n  @{ Row = 13; Col = 4; Text = \"81\u00f79=9, 0\" },\n  @{ Row = 13; Col = 5; Text = \"36\u00f78=4, 4\" },\n  @{ Row = 17; Col = 2; Text = \"58\u00f74=14, 2\" },\n  @{ Row = 17; Col = 3; Text = \"60\u00f75=12, 0\" },\n  @{ Row = 17; Col = 4; Text = \"73\u00f79=8, 1\" },\n  @{ Row = 17; Col = 5; Text = \"38\u00f78=4, 6\" }\n)\n\nforeach ($update in $cellUpdates) {\n    $cell = $t.Cell($update.Row, $update.Col)\n    $cell.Range.Text = $update.Text\n}\n"}
